$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44680
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = "`$/caja 15 kilos granel"
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 1000
$ws.Range("T2").Value = 15

# Row 3
$ws.Range("D3").Value = 44692
$ws.Range("M3").Value = 150
$ws.Range("N3").Value = 17000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 17000
$ws.Range("S3").Value = 944

# Row 4
$ws.Range("D4").Value = 44354
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100

# Row 5
$ws.Range("D5").Value = 44291
$ws.Range("M5").Value = 150
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("Q5").Value = "`$/caja 15 kilos granel"
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 800
$ws.Range("T5").Value = 15

# Row 6
$ws.Range("D6").Value = 44299
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("Q6").Value = "`$/caja 15 kilos granel"
$ws.Range("R6").Value = "Provincia de Curicó"
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 15

# Row 8
$ws.Range("D8").Value = 45062
$ws.Range("M8").Value = 200

# Row 9
$ws.Range("D9").Value = 44714
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 1111

# Row 10
$ws.Range("D10").Value = 45054
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 14000
$ws.Range("P10").Value = 14000
$ws.Range("R10").Value = "Provincia de Curicó"
$ws.Range("S10").Value = 778

# Row 11
$ws.Range("D11").Value = 44355
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 18000
$ws.Range("O11").Value = 18000
$ws.Range("P11").Value = 18000
$ws.Range("S11").Value = 1000

# Row 12
$ws.Range("D12").Value = 44326
$ws.Range("L12").Value = "Especial"
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 20000
$ws.Range("Q12").Value = "`$/caja 18 kilos granel"
$ws.Range("S12").Value = 1111
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 44340
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 230

# Row 14
$ws.Range("D14").Value = 44319
$ws.Range("M14").Value = 120
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 20000
$ws.Range("S14").Value = 1111

# Row 15
$ws.Range("D15").Value = 44348
$ws.Range("L15").Value = "Especial"
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 20000
$ws.Range("Q15").Value = "`$/caja 18 kilos granel"
$ws.Range("R15").Value = "Provincia de Limarí"
$ws.Range("S15").Value = 1111
$ws.Range("T15").Value = 18

# Row 16
$ws.Range("D16").Value = 44358
$ws.Range("L16").Value = "Especial"
$ws.Range("M16").Value = 150
$ws.Range("N16").Value = 18000
$ws.Range("O16").Value = 18000
$ws.Range("P16").Value = 18000
$ws.Range("Q16").Value = "`$/caja 18 kilos granel"
$ws.Range("R16").Value = "Provincia de Limarí"
$ws.Range("T16").Value = 18

# Row 17
$ws.Range("D17").Value = 44358
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 17000
$ws.Range("O17").Value = 17000
$ws.Range("P17").Value = 17000
$ws.Range("R17").Value = "Provincia de Limarí"
$ws.Range("S17").Value = 944

# Row 18
$ws.Range("D18").Value = 44316
$ws.Range("M18").Value = 300

# Row 19
$ws.Range("D19").Value = 44328
$ws.Range("M19").Value = 250

# Row 20
$ws.Range("D20").Value = 44691
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 17000
$ws.Range("O20").Value = 17000
$ws.Range("P20").Value = 17000
$ws.Range("S20").Value = 944

# Row 21
$ws.Range("D21").Value = 44342
$ws.Range("L21").Value = "Especial"
$ws.Range("M21").Value = 300

# Row 22
$ws.Range("D22").Value = 45083
$ws.Range("L22").Value = "Primera"
$ws.Range("N22").Value = 17000
$ws.Range("O22").Value = 17000
$ws.Range("P22").Value = 17000
$ws.Range("R22").Value = "Región de O'Higgins"
$ws.Range("S22").Value = 944
